$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.117.89"
$ws.Range("D3").Value = "'1.550.71"
$ws.Range("E3").Value = "'  -1.33%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("D6").Value = "'287.17"
$ws.Range("E6").Value = "'  -0.37%  "
$ws.Range("D7").Value = "'0.3810"
$ws.Range("E7").Value = "'  +2.44%  "
$ws.Range("D8").Value = "'0.3269"
$ws.Range("E8").Value = "'  -1.51%  "
$ws.Range("D9").Value = "'43.26"
$ws.Range("E9").Value = "'  -10.49%  "
$ws.Range("D10").Value = "'1.129"
$ws.Range("E10").Value = "'  -0.38%  "
$ws.Range("D11").Value = "'0.07326"
$ws.Range("D12").Value = "'0.9995"
$ws.Range("E12").Value = "'  -0.16%  "
$ws.Range("D13").Value = "'19.98"
$ws.Range("E13").Value = "'  -3.67%  "
$ws.Range("D14").Value = "'5.791"
$ws.Range("E14").Value = "'  -2.42%  "
$ws.Range("D15").Value = "'6.763"
$ws.Range("E15").Value = "'  -1.95%  "
$ws.Range("D16").Value = "'1.559.61"
$ws.Range("E16").Value = "'  -0.81%  "
$ws.Range("D17").Value = "'0.00001085"
$ws.Range("E17").Value = "'  -3.05%  "
$ws.Range("D18").Value = "'0.06588"
$ws.Range("E18").Value = "'  -2.31%  "
$ws.Range("D19").Value = "'85.50"
$ws.Range("E19").Value = "'  -2.65%  "
$ws.Range("D20").Value = "'0.9999"
$ws.Range("E20").Value = "'  -0.16%  "
$ws.Range("D21").Value = "'6.382"
$ws.Range("E21").Value = "'  +0.43%  "
$ws.Range("D22").Value = "'16.03"
$ws.Range("E22").Value = "'  -3.11%  "
$ws.Range("D23").Value = "'11.65"
$ws.Range("E23").Value = "'  -3.45%  "
$ws.Range("D24").Value = "'22.107.59"
$ws.Range("E24").Value = "'  -1.58%  "
$ws.Range("D25").Value = "'2.293"
$ws.Range("E25").Value = "'  -3.77%  "
$ws.Range("D26").Value = "'2.521"
$ws.Range("E26").Value = "'  -2.15%  "
$ws.Range("D27").Value = "'149.64"
$ws.Range("E27").Value = "'  -2.31%  "
$ws.Range("D28").Value = "'19.05"
$ws.Range("E28").Value = "'  -3.27%  "
$ws.Range("D30").Value = "'120.88"
$ws.Range("D31").Value = "'1.730.02"
$ws.Range("E31").Value = "'  -0.87%  "
$ws.Range("D32").Value = "'1.066"
$ws.Range("E32").Value = "'  +1.12%  "
$ws.Range("D33").Value = "'5.864"
$ws.Range("E33").Value = "'  -4.49%  "
$ws.Range("D34").Value = "'1.862"
$ws.Range("E34").Value = "'  -7.35%  "
$ws.Range("D35").Value = "'0.08212"
$ws.Range("E35").Value = "'  -1.47%  "
$ws.Range("D36").Value = "'9.272"
$ws.Range("E36").Value = "'  -5.29%  "
$ws.Range("D37").Value = "'0.02309"
$ws.Range("E37").Value = "'  -6.18%  "
$ws.Range("D38").Value = "'0.06214"
$ws.Range("E38").Value = "'  -3.05%  "
$ws.Range("D39").Value = "'5.250"
$ws.Range("E39").Value = "'  -2.00%  "
$ws.Range("E40").Value = "'  -5.41%  "
$ws.Range("D41").Value = "'1.240"
$ws.Range("E41").Value = "'  -3.89%  "
$ws.Range("D42").Value = "'10.97"
$ws.Range("E42").Value = "'  -2.90%  "
$ws.Range("E43").Value = "'  -0.09%  "
$ws.Range("D44").Value = "'0.6010"
$ws.Range("E44").Value = "'  -4.77%  "
$ws.Range("D45").Value = "'13.59"
$ws.Range("E45").Value = "'  -2.22%  "
$ws.Range("D46").Value = "'3.724"
$ws.Range("E46").Value = "'  -1.30%  "
$ws.Range("D47").Value = "'0.5796"
$ws.Range("E47").Value = "'  -5.76%  "
$ws.Range("D48").Value = "'1.972"
$ws.Range("E48").Value = "'  -4.17%  "
$ws.Range("D49").Value = "'121.77"
$ws.Range("E50").Value = "'  -3.34%  "
$ws.Range("D51").Value = "'0.07011"
$ws.Range("E51").Value = "'  -2.90%  "
